# Add a new "purchaseOrderNumber" column (H) to the payment properties sheet,
# mirroring the existing "dealerName" column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header
$ws.Range("H1").Value = "purchaseOrderNumber"

# Data rows 2..14 -> purchaseOrderNumber1 .. purchaseOrderNumber13
for ($row = 2; $row -le 14; $row++) {
    $n = $row - 1
    $ws.Cells.Item($row, 8).Value = "purchaseOrderNumber$n"
}

# Update the active selection to reflect the new column, matching the diff.
$ws.Range("H6").Select()
